# Fix for the "empty ppt" glitch in the summary pipeline:
# remove the (broken/duplicated) asset-allocation tables from both slides,
# leaving just the slide titles behind.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Delete()
        }
    }
}
